$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 0.22140880831450005
    2  = -0.0059999999672903925
    3  = -0.0039999999715742973
    4  = -0.0079999999478523875
    5  = -0.0029999999711796121
    6  = -0.0019999999697866144
    7  = -0.0099999999276403351
    8  = -0.0099999999253683747
    9  = -0.0019999999648301348
    10 = -0.0019999999626456599
    11 = -0.0029999999570522462
    12 = -0.003499999953761268
    13 = -0.0034999999515878955
    14 = -0.0079999999276134659
    15 = -0.00099999996394561919
    16 = -0.0019999999583908412
    17 = -0.0019999999580573302
    18 = -0.0039999999474504833
    19 = -0.0039999999772040162
    20 = -0.0039999999705617739
    21 = -0.042889402083706507
    22 = -0.0039999999688768995
    23 = -0.0049999999647400983
    24 = -0.019999999883452801
    25 = -0.019999999881997965
    26 = -0.0090005019581820278
    27 = -0.0024999999660009209
    28 = -0.0019999999661122203
    29 = -0.0069999999382490685
    30 = -0.059999999660635961
    31 = -0.0069999999410033098
    32 = 0.052578075282864845
    33 = -0.003999999957843059
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
